$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 1283.133333333333, 1440, 1223, 0.02807377179463704),
    @(1, 1327.433333333333, 1436, 1272, 0.02813696066538493),
    @(2, 1304.6, 1451, 1224, 0.03059898217519124),
    @(3, 1313.033333333333, 1438, 1272, 0.02815361817677816),
    @(4, 1285.366666666667, 1389, 1221, 0.0281095822652181),
    @(5, 1178, 1316, 1068, 0.02815911769866943),
    @(6, 1147.9, 1365, 1057, 0.02835159301757812),
    @(7, 1131.466666666667, 1220, 1068, 0.02830394903818766),
    @(8, 1260.9, 1375, 1221, 0.02834920088450114),
    @(9, 1152.9, 1390, 1057, 0.02836179733276367)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
